$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing rows
$ws.Range("B2").Value = "DavidAndresNorenaGiraldo"
$ws.Range("B3").Value = "LucasArboledaBedoya"

$names = @(
    "SebastianSalasCuartas",
    "JuanPabloRojasArismendy",
    "SimonLopezPelaez",
    "MiguelAngelCastillaBallestas",
    "JuanDavidRiveraCasallas",
    "EmmanuelMunozZapata",
    "JuanDavidRuizAlzate",
    "OrlandoRobertoVilladiegoOtero"
)

$row = 4
$num = 2
foreach ($name in $names) {
    # Copy the style (bold, bordered, centered) from A2 onto the new row's A cell
    $ws.Range("A2").Copy()
    $ws.Range("A" + $row).PasteSpecial(-4122)

    $ws.Range("A" + $row).Value = $num
    $ws.Range("B" + $row).Value = $name
    $ws.Range("C" + $row).Value = "No vino"
    $row++
    $num++
}
